$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 41.428665
$ws.Range("H2").Value = 124.285995
$ws.Range("I2").Value = 0.06969137269740189
$ws.Range("J2").Value = 0.06969137269740189
$ws.Range("M2").Value = 2.033666
$ws.Range("N2").Value = 6.100998000000001
$ws.Range("O2").Value = 0.01314611753252202
$ws.Range("P2").Value = 0.01314611753252202
$ws.Range("Q2").Value = 84.25206743589
$ws.Range("R2").Value = 758.26860692301
$ws.Range("S2").Value = 0.0009161709764828411
$ws.Range("T2").Value = 0.0009161709764828411
# Row 3
$ws.Range("G3").Value = 41.428665
$ws.Range("H3").Value = 124.285995
$ws.Range("I3").Value = 0.06969137269740189
$ws.Range("J3").Value = 0.06969137269740189
$ws.Range("O3").Value = 0.02278565951215614
$ws.Range("P3").Value = 0.02278565951215614
$ws.Range("Q3").Value = 146.030865541875
$ws.Range("R3").Value = 1314.277789876875
$ws.Range("S3").Value = 0.001587963889217774
$ws.Range("T3").Value = 0.001587963889217774
# Row 4
$ws.Range("G4").Value = 41.428665
$ws.Range("H4").Value = 124.285995
$ws.Range("I4").Value = 0.06969137269740189
$ws.Range("J4").Value = 0.06969137269740189
$ws.Range("M4").Value = 82.164378
$ws.Range("N4").Value = 246.493134
$ws.Range("O4").Value = 0.5311307609875792
$ws.Range("P4").Value = 0.5311307609875792
$ws.Range("Q4").Value = 3403.960491095369
$ws.Range("R4").Value = 30635.64441985833
$ws.Range("S4").Value = 0.03701523181504006
$ws.Range("T4").Value = 0.03701523181504006
# Row 5
$ws.Range("G5").Value = 41.428665
$ws.Range("H5").Value = 124.285995
$ws.Range("I5").Value = 0.06969137269740189
$ws.Range("J5").Value = 0.06969137269740189
$ws.Range("M5").Value = 1.586462
$ws.Range("N5").Value = 4.759386
$ws.Range("O5").Value = 0.01025528081448967
$ws.Range("P5").Value = 0.01025528081448967
$ws.Range("Q5").Value = 65.72500273323
$ws.Range("R5").Value = 591.52502459907
$ws.Range("S5").Value = 0.0007147045973591146
$ws.Range("T5").Value = 0.0007147045973591145
# Row 6
$ws.Range("G6").Value = 41.428665
$ws.Range("H6").Value = 124.285995
$ws.Range("I6").Value = 0.06969137269740189
$ws.Range("J6").Value = 0.06969137269740189
$ws.Range("M6").Value = 65.38769933333333
$ws.Range("N6").Value = 196.163098
$ws.Range("O6").Value = 0.422682181153253
$ws.Range("P6").Value = 0.422682181153253
$ws.Range("Q6").Value = 2708.92509080139
$ws.Range("R6").Value = 24380.32581721251
$ws.Range("S6").Value = 0.02945730141930209
$ws.Range("T6").Value = 0.02945730141930209
# Row 7
$ws.Range("H7").Value = 510.696747
$ws.Range("I7").Value = 0.2863649869040173
$ws.Range("J7").Value = 0.2863649869040173
$ws.Range("M7").Value = 2.033666
$ws.Range("N7").Value = 6.100998000000001
$ws.Range("O7").Value = 0.01314611753252202
$ws.Range("P7").Value = 0.01314611753252202
$ws.Range("Q7").Value = 346.195536894834
$ws.Range("R7").Value = 3115.759832053506
$ws.Range("S7").Value = 0.003764587775039339
$ws.Range("T7").Value = 0.00376458777503934
# Row 8
$ws.Range("H8").Value = 510.696747
$ws.Range("I8").Value = 0.2863649869040173
$ws.Range("J8").Value = 0.2863649869040173
$ws.Range("O8").Value = 0.02278565951215614
$ws.Range("P8").Value = 0.02278565951215614
$ws.Range("Q8").Value = 600.0473986938749
$ws.Range("R8").Value = 5400.426588244874
$ws.Range("S8").Value = 0.006525015087797991
$ws.Range("T8").Value = 0.006525015087797992
# Row 9
$ws.Range("H9").Value = 510.696747
$ws.Range("I9").Value = 0.2863649869040173
$ws.Range("J9").Value = 0.2863649869040173
$ws.Range("M9").Value = 82.164378
$ws.Range("N9").Value = 246.493134
$ws.Range("O9").Value = 0.5311307609875792
$ws.Range("P9").Value = 0.5311307609875792
$ws.Range("Q9").Value = 13987.02685462612
$ws.Range("R9").Value = 125883.2416916351
$ws.Range("S9").Value = 0.1520972534145288
$ws.Range("T9").Value = 0.1520972534145289
# Row 10
$ws.Range("H10").Value = 510.696747
$ws.Range("I10").Value = 0.2863649869040173
$ws.Range("J10").Value = 0.2863649869040173
$ws.Range("M10").Value = 1.586462
$ws.Range("N10").Value = 4.759386
$ws.Range("O10").Value = 0.01025528081448967
$ws.Range("P10").Value = 0.01025528081448967
$ws.Range("Q10").Value = 270.066994213038
$ws.Range("R10").Value = 2430.602947917342
$ws.Range("S10").Value = 0.002936753356138353
$ws.Range("T10").Value = 0.002936753356138353
# Row 11
$ws.Range("H11").Value = 510.696747
$ws.Range("I11").Value = 0.2863649869040173
$ws.Range("J11").Value = 0.2863649869040173
$ws.Range("M11").Value = 65.38769933333333
$ws.Range("N11").Value = 196.163098
$ws.Range("O11").Value = 0.422682181153253
$ws.Range("P11").Value = 0.422682181153253
$ws.Range("Q11").Value = 11131.09511444913
$ws.Range("R11").Value = 100179.8560300422
$ws.Range("S11").Value = 0.1210413772705128
$ws.Range("T11").Value = 0.1210413772705128
# Row 12
$ws.Range("G12").Value = 244.5761666666666
$ws.Range("H12").Value = 733.7284999999999
$ws.Range("I12").Value = 0.4114264551867299
$ws.Range("J12").Value = 0.41142645518673
$ws.Range("M12").Value = 2.033666
$ws.Range("N12").Value = 6.100998000000001
$ws.Range("O12").Value = 0.01314611753252202
$ws.Range("P12").Value = 0.01314611753252202
$ws.Range("Q12").Value = 497.3862345603333
$ws.Range("R12").Value = 4476.476111043
$ws.Range("S12").Value = 0.005408660535873655
$ws.Range("T12").Value = 0.005408660535873655
# Row 13
$ws.Range("G13").Value = 244.5761666666666
$ws.Range("H13").Value = 733.7284999999999
$ws.Range("I13").Value = 0.4114264551867299
$ws.Range("J13").Value = 0.41142645518673
$ws.Range("O13").Value = 0.02278565951215614
$ws.Range("P13").Value = 0.02278565951215614
$ws.Range("Q13").Value = 862.1004154791665
$ws.Range("R13").Value = 7758.903739312499
$ws.Range("S13").Value = 0.009374623122178197
$ws.Range("T13").Value = 0.009374623122178197
# Row 14
$ws.Range("G14").Value = 244.5761666666666
$ws.Range("H14").Value = 733.7284999999999
$ws.Range("I14").Value = 0.4114264551867299
$ws.Range("J14").Value = 0.41142645518673
$ws.Range("M14").Value = 82.164378
$ws.Range("N14").Value = 246.493134
$ws.Range("O14").Value = 0.5311307609875792
$ws.Range("P14").Value = 0.5311307609875792
$ws.Range("Q14").Value = 20095.448607791
$ws.Range("R14").Value = 180859.037470119
$ws.Range("S14").Value = 0.21852124623375
$ws.Range("T14").Value = 0.21852124623375
# Row 15
$ws.Range("G15").Value = 244.5761666666666
$ws.Range("H15").Value = 733.7284999999999
$ws.Range("I15").Value = 0.4114264551867299
$ws.Range("J15").Value = 0.41142645518673
$ws.Range("M15").Value = 1.586462
$ws.Range("N15").Value = 4.759386
$ws.Range("O15").Value = 0.01025528081448967
$ws.Range("P15").Value = 0.01025528081448967
$ws.Range("Q15").Value = 388.0107945223333
$ws.Range("R15").Value = 3492.097150701
$ws.Range("S15").Value = 0.004219293832449964
$ws.Range("T15").Value = 0.004219293832449964
# Row 16
$ws.Range("G16").Value = 244.5761666666666
$ws.Range("H16").Value = 733.7284999999999
$ws.Range("I16").Value = 0.4114264551867299
$ws.Range("J16").Value = 0.41142645518673
$ws.Range("M16").Value = 65.38769933333333
$ws.Range("N16").Value = 196.163098
$ws.Range("O16").Value = 0.422682181153253
$ws.Range("P16").Value = 0.422682181153253
$ws.Range("Q16").Value = 15992.27285009922
$ws.Range("R16").Value = 143930.455650893
$ws.Range("S16").Value = 0.1739026314624781
$ws.Range("T16").Value = 0.1739026314624781
# Row 17
$ws.Range("G17").Value = 24.173247
$ws.Range("H17").Value = 72.51974100000001
$ws.Range("I17").Value = 0.04066427836821081
$ws.Range("J17").Value = 0.04066427836821081
$ws.Range("M17").Value = 2.033666
$ws.Range("N17").Value = 6.100998000000001
$ws.Range("O17").Value = 0.01314611753252202
$ws.Range("P17").Value = 0.01314611753252202
$ws.Range("Q17").Value = 49.16031053350201
$ws.Range("R17").Value = 442.4427948015181
$ws.Range("S17").Value = 0.000534577382803692
$ws.Range("T17").Value = 0.000534577382803692
# Row 18
$ws.Range("G18").Value = 24.173247
$ws.Range("H18").Value = 72.51974100000001
$ws.Range("I18").Value = 0.04066427836821081
$ws.Range("J18").Value = 0.04066427836821081
$ws.Range("O18").Value = 0.02278565951215614
$ws.Range("P18").Value = 0.02278565951215614
$ws.Range("Q18").Value = 85.20767401912501
$ws.Range("R18").Value = 766.869066172125
$ws.Range("S18").Value = 0.0009265624012055881
$ws.Range("T18").Value = 0.0009265624012055881
# Row 19
$ws.Range("G19").Value = 24.173247
$ws.Range("H19").Value = 72.51974100000001
$ws.Range("I19").Value = 0.04066427836821081
$ws.Range("J19").Value = 0.04066427836821081
$ws.Range("M19").Value = 82.164378
$ws.Range("N19").Value = 246.493134
$ws.Range("O19").Value = 0.5311307609875792
$ws.Range("P19").Value = 0.5311307609875792
$ws.Range("Q19").Value = 1986.179803995366
$ws.Range("R19").Value = 17875.6182359583
$ws.Range("S19").Value = 0.02159804911471856
$ws.Range("T19").Value = 0.02159804911471856
# Row 20
$ws.Range("G20").Value = 24.173247
$ws.Range("H20").Value = 72.51974100000001
$ws.Range("I20").Value = 0.04066427836821081
$ws.Range("J20").Value = 0.04066427836821081
$ws.Range("M20").Value = 1.586462
$ws.Range("N20").Value = 4.759386
$ws.Range("O20").Value = 0.01025528081448967
$ws.Range("P20").Value = 0.01025528081448967
$ws.Range("Q20").Value = 38.349937782114
$ws.Range("R20").Value = 345.149440039026
$ws.Range("S20").Value = 0.0004170235937845795
$ws.Range("T20").Value = 0.0004170235937845795
# Row 21
$ws.Range("G21").Value = 24.173247
$ws.Range("H21").Value = 72.51974100000001
$ws.Range("I21").Value = 0.04066427836821081
$ws.Range("J21").Value = 0.04066427836821081
$ws.Range("M21").Value = 65.38769933333333
$ws.Range("N21").Value = 196.163098
$ws.Range("O21").Value = 0.422682181153253
$ws.Range("P21").Value = 0.422682181153253
$ws.Range("Q21").Value = 1580.633006746402
$ws.Range("R21").Value = 14225.69706071762
$ws.Range("S21").Value = 0.01718806587569839
$ws.Range("T21").Value = 0.01718806587569839
# Row 22
$ws.Range("G22").Value = 114.0486906666667
$ws.Range("H22").Value = 342.146072
$ws.Range("I22").Value = 0.19185290684364
$ws.Range("J22").Value = 0.19185290684364
$ws.Range("M22").Value = 2.033666
$ws.Range("N22").Value = 6.100998000000001
$ws.Range("O22").Value = 0.01314611753252202
$ws.Range("P22").Value = 0.01314611753252202
$ws.Range("Q22").Value = 231.9369445533174
$ws.Range("R22").Value = 2087.432500979856
$ws.Range("S22").Value = 0.002522120862322489
$ws.Range("T22").Value = 0.002522120862322489
# Row 23
$ws.Range("G23").Value = 114.0486906666667
$ws.Range("H23").Value = 342.146072
$ws.Range("I23").Value = 0.19185290684364
$ws.Range("J23").Value = 0.19185290684364
$ws.Range("O23").Value = 0.02278565951215614
$ws.Range("P23").Value = 0.02278565951215614
$ws.Range("Q23").Value = 402.0073785136667
$ws.Range("R23").Value = 3618.066406623
$ws.Range("S23").Value = 0.004371495011756592
$ws.Range("T23").Value = 0.004371495011756592
# Row 24
$ws.Range("G24").Value = 114.0486906666667
$ws.Range("H24").Value = 342.146072
$ws.Range("I24").Value = 0.19185290684364
$ws.Range("J24").Value = 0.19185290684364
$ws.Range("M24").Value = 82.164378
$ws.Range("N24").Value = 246.493134
$ws.Range("O24").Value = 0.5311307609875792
$ws.Range("P24").Value = 0.5311307609875792
$ws.Range("Q24").Value = 9370.739730341073
$ws.Range("R24").Value = 84336.65757306965
$ws.Range("S24").Value = 0.1018989804095416
$ws.Range("T24").Value = 0.1018989804095416
# Row 25
$ws.Range("G25").Value = 114.0486906666667
$ws.Range("H25").Value = 342.146072
$ws.Range("I25").Value = 0.19185290684364
$ws.Range("J25").Value = 0.19185290684364
$ws.Range("M25").Value = 1.586462
$ws.Range("N25").Value = 4.759386
$ws.Range("O25").Value = 0.01025528081448967
$ws.Range("P25").Value = 0.01025528081448967
$ws.Range("Q25").Value = 180.9339138924213
$ws.Range("R25").Value = 1628.405225031792
$ws.Range("S25").Value = 0.001967505434757654
$ws.Range("T25").Value = 0.001967505434757654
# Row 26
$ws.Range("G26").Value = 114.0486906666667
$ws.Range("H26").Value = 342.146072
$ws.Range("I26").Value = 0.19185290684364
$ws.Range("J26").Value = 0.19185290684364
$ws.Range("M26").Value = 65.38769933333333
$ws.Range("N26").Value = 196.163098
$ws.Range("O26").Value = 0.422682181153253
$ws.Range("P26").Value = 0.422682181153253
$ws.Range("Q26").Value = 7457.38149467234
$ws.Range("R26").Value = 67116.43345205106
$ws.Range("S26").Value = 0.08109280512526161
$ws.Range("T26").Value = 0.0810928051252616
